$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (P1) entirely, shifting rows 3-5 up to 2-4
$ws.Rows.Item(2).Delete()

# Edit resulting row 2 (was P2 row): bump "Pelayanan" value
$ws.Range("E2").Value = 100

# Edit resulting row 3 (previously blank-name row) into P3
$ws.Range("A3").Value = "P3"
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 10

# Edit resulting row 4 (previously blank-name row) into P4
$ws.Range("A4").Value = "P4"
$ws.Range("E4").Value = 100
